$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.028139083075631
$ws.Range("D2").Value = 1.029433777072122
$ws.Range("E2").Value = 1.03741019972774
$ws.Range("F2").Value = 1.04664831816462
$ws.Range("I2").Value = 1.032572425840893
$ws.Range("J2").Value = 1.033293061461415
$ws.Range("K2").Value = 1.032247700814112
$ws.Range("L2").Value = 1.040201163811946
$ws.Range("M2").Value = 1.049413165436673
$ws.Range("N2").Value = 1.015107105878073

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.028978084316564
$ws.Range("D3").Value = 1.030147519223694
$ws.Range("E3").Value = 1.038163285338898
$ws.Range("F3").Value = 1.047502073042091
$ws.Range("I3").Value = 1.032662048692641
$ws.Range("J3").Value = 1.033773111755603
$ws.Range("K3").Value = 1.032769841641437
$ws.Range("L3").Value = 1.040764201503886
$ws.Range("M3").Value = 1.050078495992177
$ws.Range("N3").Value = 1.015267485934195

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.029521606032618
$ws.Range("D4").Value = 1.030610230720443
$ws.Range("E4").Value = 1.038651530383783
$ws.Range("F4").Value = 1.048055547154961
$ws.Range("I4").Value = 1.032718755058449
$ws.Range("J4").Value = 1.034083717707756
$ws.Range("K4").Value = 1.0331079024437
$ws.Range("L4").Value = 1.041128808075011
$ws.Range("M4").Value = 1.050509402095662
$ws.Range("N4").Value = 1.015371215424434

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.02975025215044
$ws.Range("D5").Value = 1.030804961687856
$ws.Range("E5").Value = 1.038857013742746
$ws.Range("F5").Value = 1.048288474078009
$ws.Range("I5").Value = 1.032742285980159
$ws.Range("J5").Value = 1.034214290657588
$ws.Range("K5").Value = 1.033250069851291
$ws.Range("L5").Value = 1.041282155068073
$ws.Range("M5").Value = 1.050690647337268
$ws.Range("N5").Value = 1.015414811536827

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.029788651563202
$ws.Range("D6").Value = 1.030837669962144
$ws.Range("E6").Value = 1.038891528450577
$ws.Range("F6").Value = 1.048327597924475
$ws.Range("I6").Value = 1.03274621881964
$ws.Range("J6").Value = 1.034236214052402
$ws.Range("K6").Value = 1.033273943078681
$ws.Range("L6").Value = 1.041307906552857
$ws.Range("M6").Value = 1.050721084574618
$ws.Range("N6").Value = 1.015422130811978

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.029524660627541
$ws.Range("D7").Value = 1.030612831913903
$ws.Range("E7").Value = 1.038654275180195
$ws.Range("F7").Value = 1.048058658570215
$ws.Range("I7").Value = 1.032719070692532
$ws.Range("J7").Value = 1.034085462454319
$ws.Range("K7").Value = 1.033109801910054
$ws.Range("L7").Value = 1.041130856846038
$ws.Range("M7").Value = 1.050511823542859
$ws.Range("N7").Value = 1.015371798004657

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.02842249564873
$ws.Range("D8").Value = 1.029674807786493
$ws.Range("E8").Value = 1.037664511053574
$ws.Range("F8").Value = 1.046936632573459
$ws.Range("I8").Value = 1.032602980103009
$ws.Range("J8").Value = 1.033455299782875
$ws.Range("K8").Value = 1.032424118327215
$ws.Range("L8").Value = 1.040391385433597
$ws.Range("M8").Value = 1.049637934665668
$ws.Range("N8").Value = 1.015161316537826

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.026485252902074
$ws.Range("D9").Value = 1.028028654143275
$ws.Range("E9").Value = 1.035927766232374
$ws.Range("F9").Value = 1.044967514866062
$ws.Range("I9").Value = 1.032388596429427
$ws.Range("J9").Value = 1.032344787797795
$ws.Range("K9").Value = 1.031217456645714
$ws.Range("L9").Value = 1.039090579663915
$ws.Range("M9").Value = 1.048101114301258
$ws.Range("N9").Value = 1.014790081642496

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.025197160181891
$ws.Range("D10").Value = 1.026935876838196
$ws.Range("E10").Value = 1.034774987725748
$ws.Range("F10").Value = 1.043660291782816
$ws.Range("I10").Value = 1.032239112582099
$ws.Range("J10").Value = 1.03160446802138
$ws.Range("K10").Value = 1.030414176154306
$ws.Range("L10").Value = 1.038224968841821
$ws.Range("M10").Value = 1.047078742461411
$ws.Range("N10").Value = 1.014542392076196

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.024640229724025
$ws.Range("D11").Value = 1.026463820435856
$ws.Range("E11").Value = 1.034277043769476
$ws.Range("F11").Value = 1.043095584750406
$ws.Range("I11").Value = 1.03217283639085
$ws.Range("J11").Value = 1.03128392307303
$ws.Range("K11").Value = 1.030066639982275
$ws.Range("L11").Value = 1.037850546384501
$ws.Range("M11").Value = 1.046636581581293
$ws.Range("N11").Value = 1.014435098470541

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.024433486155525
$ws.Range("D12").Value = 1.026288648338371
$ws.Range("E12").Value = 1.034092269943074
$ws.Range("F12").Value = 1.042886029155947
$ws.Range("I12").Value = 1.03214798643459
$ws.Range("J12").Value = 1.031164862504892
$ws.Range("K12").Value = 1.029937594522036
$ws.Range("L12").Value = 1.037711529730735
$ws.Range("M12").Value = 1.046472425138889
$ws.Range("N12").Value = 1.014395238967527

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.02447782766023
$ws.Range("D13").Value = 1.026326215594429
$ws.Range("E13").Value = 1.034131896165715
$ws.Range("F13").Value = 1.042930970380273
$ws.Range("I13").Value = 1.032153327323925
$ws.Range("J13").Value = 1.031190401194073
$ws.Range("K13").Value = 1.02996527315418
$ws.Range("L13").Value = 1.037741346512016
$ws.Range("M13").Value = 1.046507633523708
$ws.Range("N13").Value = 1.014403789224514

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.024623137672394
$ws.Range("D14").Value = 1.026449337163104
$ws.Range("E14").Value = 1.034261766527769
$ws.Range("F14").Value = 1.043078258692316
$ws.Range("I14").Value = 1.032170787015269
$ws.Range("J14").Value = 1.03127408139802
$ws.Range("K14").Value = 1.030055972119014
$ws.Range("L14").Value = 1.037839053981237
$ws.Range("M14").Value = 1.046623010679901
$ws.Range("N14").Value = 1.014431803788732

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.024712684607245
$ws.Range("D15").Value = 1.026525219096929
$ws.Range("E15").Value = 1.034341808480232
$ws.Range("F15").Value = 1.043169034676298
$ws.Range("I15").Value = 1.032181513780757
$ws.Range("J15").Value = 1.031325640116364
$ws.Range("K15").Value = 1.030111860743858
$ws.Range("L15").Value = 1.037899262842252
$ws.Range("M15").Value = 1.046694109244444
$ws.Range("N15").Value = 1.014449063721686

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.025234139287776
$ws.Range("D16").Value = 1.026967229496964
$ws.Range("E16").Value = 1.034808060426191
$ws.Range("F16").Value = 1.043697797707823
$ws.Range("I16").Value = 1.03224347854051
$ws.Range("J16").Value = 1.031625742035423
$ws.Range("K16").Value = 1.030437247209702
$ws.Range("L16").Value = 1.038249826430964
$ws.Range("M16").Value = 1.047108098598152
$ws.Range("N16").Value = 1.014549511948007

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.025561455338156
$ws.Range("D17").Value = 1.027244793041629
$ws.Range("E17").Value = 1.035100855072164
$ws.Range("F17").Value = 1.044029834093431
$ws.Range("I17").Value = 1.0322819332228
$ws.Range("J17").Value = 1.031813993908454
$ws.Range("K17").Value = 1.030641432090463
$ws.Range("L17").Value = 1.038469831924735
$ws.Range("M17").Value = 1.047367927074288
$ws.Range("N17").Value = 1.014612509462401

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.025752452469449
$ws.Range("D18").Value = 1.027406799424464
$ws.Range("E18").Value = 1.035271754504454
$ws.Range("F18").Value = 1.044223633409862
$ws.Range("I18").Value = 1.03230421371799
$ws.Range("J18").Value = 1.031923799726972
$ws.Range("K18").Value = 1.030760557520669
$ws.Range("L18").Value = 1.038598195199411
$ws.Range("M18").Value = 1.047519531926498
$ws.Range("N18").Value = 1.01464925072411

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.025817590941017
$ws.Range("D19").Value = 1.027462057720083
$ws.Range("E19").Value = 1.035330046626607
$ws.Range("F19").Value = 1.044289735636044
$ws.Range("I19").Value = 1.032311785418637
$ws.Range("J19").Value = 1.03196124091884
$ws.Range("K19").Value = 1.030801180892806
$ws.Range("L19").Value = 1.038641970123064
$ws.Range("M19").Value = 1.047571233911451
$ws.Range("N19").Value = 1.014661777830758

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.025526329217769
$ws.Range("D20").Value = 1.027215001888286
$ws.Range("E20").Value = 1.035069428834039
$ws.Range("F20").Value = 1.043994196465806
$ws.Range("I20").Value = 1.032277822855081
$ws.Range("J20").Value = 1.031793796080066
$ws.Range("K20").Value = 1.030619522094724
$ws.Range("L20").Value = 1.038446223497403
$ws.Range("M20").Value = 1.047340044619256
$ws.Range("N20").Value = 1.014605750848208

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.024580344025936
$ws.Range("D21").Value = 1.026413076203371
$ws.Range("E21").Value = 1.034223517828155
$ws.Range("F21").Value = 1.043034880370023
$ws.Range("I21").Value = 1.032165651971249
$ws.Range("J21").Value = 1.031249439539802
$ws.Range("K21").Value = 1.030029262282122
$ws.Range("L21").Value = 1.037810279884104
$ws.Range("M21").Value = 1.04658903271119
$ws.Range("N21").Value = 1.014423554356234

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.023986290795356
$ws.Range("D22").Value = 1.025909862068527
$ws.Range("E22").Value = 1.033692729747902
$ws.Range("F22").Value = 1.042432888488517
$ws.Range("I22").Value = 1.032093783563331
$ws.Range("J22").Value = 1.030907205346657
$ws.Range("K22").Value = 1.029658403279318
$ws.Range("L22").Value = 1.037410787781024
$ws.Range("M22").Value = 1.046117315821333
$ws.Range("N22").Value = 1.014308966407957

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.024301140351166
$ws.Range("D23").Value = 1.026176531084787
$ws.Range("E23").Value = 1.033974008509533
$ws.Range("F23").Value = 1.042751904358082
$ws.Range("I23").Value = 1.032132009367654
$ws.Range("J23").Value = 1.031088627412664
$ws.Range("K23").Value = 1.029854977454598
$ws.Range("L23").Value = 1.037622532338418
$ws.Range("M23").Value = 1.046367336350167
$ws.Range("N23").Value = 1.014369714686849

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.025542200955018
$ws.Range("D24").Value = 1.027228462893866
$ws.Range("E24").Value = 1.035083628636227
$ws.Range("F24").Value = 1.044010299180144
$ws.Range("I24").Value = 1.032279680615242
$ws.Range("J24").Value = 1.031802922604265
$ws.Range("K24").Value = 1.030629422192842
$ws.Range("L24").Value = 1.038456891012919
$ws.Range("M24").Value = 1.047352643342212
$ws.Range("N24").Value = 1.01460880478784

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.026985484170102
$ws.Range("D25").Value = 1.028453411296939
$ws.Range("E25").Value = 1.036375874276728
$ws.Range("F25").Value = 1.044967514866062
$ws.Range("I25").Value = 1.032445179107892
$ws.Range("J25").Value = 1.032631883719522
$ws.Range("K25").Value = 1.031529209236345
$ws.Range("L25").Value = 1.039426594850307
$ws.Range("M25").Value = 1.048498043291606
$ws.Range("N25").Value = 1.014886092050284
